$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.215.58'
$ws.Range('E2').Value = '  +1.09%  '

$ws.Range('D3').Value = '1.652.43'
$ws.Range('E3').Value = '  +0.90%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +1.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.87'
$ws.Range('E5').Value = '  +0.67%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.508'
$ws.Range('E6').Value = '  +0.67%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('E7').Value = '  +1.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.0638'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.254'
$ws.Range('E9').Value = '  -0.78%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.07'
$ws.Range('E10').Value = '  -2.75%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0799'
$ws.Range('E11').Value = '  +0.70%  '

$ws.Range('D12').Value = '1.717.17'
$ws.Range('E12').Value = '  +4.93%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.26'
$ws.Range('E13').Value = '  +0.00%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.541'
$ws.Range('E14').Value = '  -0.41%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '63.49'
$ws.Range('E15').Value = '  +0.86%  '

$ws.Range('D16').Value = '0.0₃0762'
$ws.Range('E16').Value = '  -0.21%  '

$ws.Range('D17').Value = '26.244.29'
$ws.Range('E17').Value = '  +1.33%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.01'
$ws.Range('E18').Value = '  +1.18%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '194.14'
$ws.Range('E19').Value = '  +0.70%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.33'
$ws.Range('E20').Value = '  -0.69%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.74'
$ws.Range('E21').Value = '  -1.76%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.18'
$ws.Range('E22').Value = '  -1.66%  '

$ws.Range('B23').Value = 'Monero'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '144.82'
$ws.Range('E23').Value = '  +0.46%  '

$ws.Range('B24').Value = 'Stellar'
$ws.Range('C24').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.131'
$ws.Range('E24').Value = '  +0.94%  '

$ws.Range('E25').Value = '  +1.55%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.78'
$ws.Range('E26').Value = '  -0.74%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.84'
$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.51'
$ws.Range('E28').Value = '  -0.24%  '

$ws.Range('E29').Value = '  +0.39%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0488'
$ws.Range('E30').Value = '  -2.72%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.26'
$ws.Range('E31').Value = '  +0.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.27'
$ws.Range('E32').Value = '  -1.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.53'
$ws.Range('E33').Value = '  +0.19%  '

$ws.Range('E34').Value = '  +1.34%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.905'
$ws.Range('E35').Value = '  +0.42%  '

$ws.Range('D36').Value = '1.138.02'
$ws.Range('E36').Value = '  +0.22%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.48'
$ws.Range('E37').Value = '  +1.02%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.531'
$ws.Range('E38').Value = '  -2.68%  '

$ws.Range('E39').Value = '  -0.78%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.800'
$ws.Range('E40').Value = '  +0.37%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.85'
$ws.Range('E41').Value = '  -0.44%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.33'
$ws.Range('E42').Value = '  -2.75%  '

$ws.Range('D43').Value = '0.0₆0115'
$ws.Range('E43').Value = '  +0.95%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '56.32'
$ws.Range('E44').Value = '  -0.54%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.51'
$ws.Range('E45').Value = '  +1.63%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0523'
$ws.Range('E46').Value = '  -1.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.79'
$ws.Range('E47').Value = '  +1.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.418'
$ws.Range('E48').Value = '  +0.85%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0938'
$ws.Range('E50').Value = '  -2.67%  '

$ws.Range('E51').Value = '  +1.51%  '
